$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 367.2
$ws.Range("J2").Value = 445.33334
$ws.Range("L2").Value = 445.33334
$ws.Range("N2").Value = -671.33334
$ws.Range("H40").Value = 2964.625
$ws.Range("I40").Value = 2913.3333
$ws.Range("J40").Value = 2995.4
$ws.Range("K40").Value = 2913.3333
$ws.Range("L40").Value = 2995.4
$ws.Range("M40").Value = -2738.3333
$ws.Range("N40").Value = -3345.4
$ws.Range("H43").Value = 1421.4286
$ws.Range("J43").Value = 1375
$ws.Range("L43").Value = 1375
$ws.Range("N43").Value = -1513
$ws.Range("H53").Value = 5099.25
$ws.Range("I53").Value = 5099.25
$ws.Range("K53").Value = 5099.25
$ws.Range("M53").Value = -4462.25
$ws.Range("H62").Value = 24419.5
$ws.Range("I62").Value = 24419.5
$ws.Range("K62").Value = 24419.5
$ws.Range("M62").Value = -23795.5
$ws.Range("H64").Value = 3462.5
$ws.Range("I64").Value = 3060
$ws.Range("K64").Value = 3060
$ws.Range("M64").Value = -2812
$ws.Range("H65").Value = 24419.5
$ws.Range("I65").Value = 24419.5
$ws.Range("K65").Value = 122097.5
$ws.Range("M65").Value = -118977.5
$ws.Range("H67").Value = 3462.5
$ws.Range("I67").Value = 3060
$ws.Range("K67").Value = 3060
$ws.Range("M67").Value = -2202
$ws.Range("H76").Value = 3556.8572
$ws.Range("I76").Value = 3419.8
$ws.Range("K76").Value = 3419.8
$ws.Range("M76").Value = -3104.8
$ws.Range("H79").Value = 3556.8572
$ws.Range("I79").Value = 3419.8
$ws.Range("K79").Value = 3419.8
$ws.Range("M79").Value = -2327.8
$ws.Range("H100").Value = 815.82355
$ws.Range("J100").Value = 816.8
$ws.Range("L100").Value = 816.8
$ws.Range("N100").Value = -1898.8
$ws.Range("H116").Value = 10175
$ws.Range("I116").Value = 18862.834
$ws.Range("J116").Value = 2728.2856
$ws.Range("K116").Value = 18862.834
$ws.Range("L116").Value = 2728.2856
$ws.Range("M116").Value = -15420.834
$ws.Range("N116").Value = -9612.285599999999
$ws.Range("H125").Value = 1066.591
$ws.Range("I125").Value = 837.5
$ws.Range("J125").Value = 1677.5
$ws.Range("K125").Value = 7537.5
$ws.Range("L125").Value = 15097.5
$ws.Range("M125").Value = -5077.5
$ws.Range("N125").Value = -20017.5
$ws.Range("H132").Value = 926.25
$ws.Range("I132").Value = 768
$ws.Range("J132").Value = 1875.75
$ws.Range("K132").Value = 2304
$ws.Range("L132").Value = 5627.25
$ws.Range("M132").Value = 226
$ws.Range("N132").Value = -10687.25
$ws.Range("H138").Value = 1515.04
$ws.Range("I138").Value = 1248.4517
$ws.Range("J138").Value = 1634.8116
$ws.Range("K138").Value = 3745.3551
$ws.Range("L138").Value = 4904.4348
$ws.Range("M138").Value = 1394.6449
$ws.Range("N138").Value = -15184.4348

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5536.22
$ws.Range("I32").Value = 3541.4111
$ws.Range("J32").Value = 23489.5
$ws.Range("K32").Value = 3541.4111
$ws.Range("L32").Value = 23489.5
$ws.Range("M32").Value = -3254.4111
$ws.Range("N32").Value = -24063.5
$ws.Range("H45").Value = 1632.3334
$ws.Range("I45").Value = 1259.2
$ws.Range("K45").Value = 1259.2
$ws.Range("M45").Value = -882.2
$ws.Range("H61").Value = 28653.934
$ws.Range("I61").Value = 34212.875
$ws.Range("J61").Value = 6418.1665
$ws.Range("K61").Value = 34212.875
$ws.Range("L61").Value = 6418.1665
$ws.Range("M61").Value = -34000.875
$ws.Range("N61").Value = -6842.1665
$ws.Range("H122").Value = 2810
$ws.Range("I122").Value = 1503.6666
$ws.Range("J122").Value = 5749.25
$ws.Range("K122").Value = 4510.9998
$ws.Range("L122").Value = 17247.75
$ws.Range("M122").Value = -2060.9998
$ws.Range("N122").Value = -22147.75
$ws.Range("H136").Value = 28653.934
$ws.Range("I136").Value = 34212.875
$ws.Range("J136").Value = 6418.1665
$ws.Range("K136").Value = 102638.625
$ws.Range("L136").Value = 19254.4995
$ws.Range("M136").Value = -100088.625
$ws.Range("N136").Value = -24354.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6362.8716
$ws.Range("I134").Value = 7013.448
$ws.Range("K134").Value = 21040.344
$ws.Range("M134").Value = -18505.344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 509.23077
$ws.Range("I105").Value = 509.23077
$ws.Range("K105").Value = 509.23077
$ws.Range("M105").Value = 1237.76923
$ws.Range("H132").Value = 1418.4828
$ws.Range("I132").Value = 853.5599999999999
$ws.Range("K132").Value = 2560.68
$ws.Range("M132").Value = -30.67999999999984
$ws.Range("H134").Value = 1259.9269
$ws.Range("I134").Value = 1134.8485
$ws.Range("K134").Value = 3404.5455
$ws.Range("M134").Value = -869.5455000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2744.8572
$ws.Range("J81").Value = 3002.3333
$ws.Range("L81").Value = 9006.999899999999
$ws.Range("N81").Value = -11252.9999
$ws.Range("H84").Value = 2744.8572
$ws.Range("J84").Value = 3002.3333
$ws.Range("L84").Value = 27020.9997
$ws.Range("N84").Value = -38252.9997
$ws.Range("H92").Value = 532.6667
$ws.Range("J92").Value = 600
$ws.Range("L92").Value = 1800
$ws.Range("N92").Value = -4296
$ws.Range("H113").Value = 56118.25
$ws.Range("I113").Value = 220530.8
$ws.Range("J113").Value = 1314.0667
$ws.Range("K113").Value = 661592.3999999999
$ws.Range("L113").Value = 3942.2001
$ws.Range("M113").Value = -659422.3999999999
$ws.Range("N113").Value = -8282.2001
$ws.Range("H131").Value = 15474.928
$ws.Range("J131").Value = 17321.06
$ws.Range("L131").Value = 51963.18000000001
$ws.Range("N131").Value = -62043.18000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1023
$ws.Range("I97").Value = 987.1177
$ws.Range("J97").Value = 1099.25
$ws.Range("K97").Value = 987.1177
$ws.Range("L97").Value = 1099.25
$ws.Range("M97").Value = -491.1177
$ws.Range("N97").Value = -2091.25
$ws.Range("H122").Value = 1387.375
$ws.Range("I122").Value = 1349.8334
$ws.Range("K122").Value = 4049.5002
$ws.Range("M122").Value = -1599.5002
$ws.Range("H132").Value = 715936.9
$ws.Range("I132").Value = 858325.1
$ws.Range("K132").Value = 2574975.3
$ws.Range("M132").Value = -2572445.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2614.9
$ws.Range("I46").Value = 1787.5
$ws.Range("K46").Value = 1787.5
$ws.Range("M46").Value = -1599.5
$ws.Range("H132").Value = 3662.842
$ws.Range("I132").Value = 1965.7142
$ws.Range("K132").Value = 5897.142599999999
$ws.Range("M132").Value = -3367.142599999999
$ws.Range("H136").Value = 2791.7646
$ws.Range("I136").Value = 2777.3333
$ws.Range("K136").Value = 8331.999899999999
$ws.Range("M136").Value = -5781.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47087
$ws.Range("J70").Value = 47087
$ws.Range("L70").Value = 47087
$ws.Range("N70").Value = -47717
$ws.Range("H73").Value = 47087
$ws.Range("J73").Value = 47087
$ws.Range("L73").Value = 47087
$ws.Range("N73").Value = -49271
$ws.Range("H107").Value = 594.5263
$ws.Range("I107").Value = 516.44446
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1549.33338
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 370.66662
$ws.Range("N107").Value = -9840
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H126").Value = 1588.6774
$ws.Range("I126").Value = 1237.0416
$ws.Range("J126").Value = 2794.2856
$ws.Range("K126").Value = 3711.1248
$ws.Range("L126").Value = 8382.856800000001
$ws.Range("M126").Value = -1241.1248
$ws.Range("N126").Value = -13322.8568
$ws.Range("H132").Value = 1065.6025
$ws.Range("I132").Value = 966.8939
$ws.Range("J132").Value = 1608.5
$ws.Range("K132").Value = 2900.6817
$ws.Range("L132").Value = 4825.5
$ws.Range("M132").Value = -370.6817000000001
$ws.Range("N132").Value = -9885.5
$ws.Range("H136").Value = 1772.68
$ws.Range("I136").Value = 1175.4117
$ws.Range("J136").Value = 3041.875
$ws.Range("K136").Value = 3526.2351
$ws.Range("L136").Value = 9125.625
$ws.Range("M136").Value = -976.2351000000003
$ws.Range("N136").Value = -14225.625
